# Apply updated cryptocurrency price/volume data to sheet1 (cells B/C/D/E, rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.627.81"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.726.72"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D7").Value = "3.726.94"
$ws.Range("E7").Value = "  -1.88%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  -4.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").Value = "4.352.14"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "3.734.73"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "67.551.23"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.49%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.698"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.51%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("E25").Value = "  -10.63%  "
$ws.Range("E26").Value = "  -5.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "3.873.11"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.84%  "
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.90%  "
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").Value = "3.678.09"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.137"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.990"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "385.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0345"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "
